$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; old C..F shift to D..G.
$ws.Columns("C").Insert()

# Header text for the new "Trạng thái" (Status) column.
$ws.Cells.Item(1, 3).Value = "Trạng thái"

# Status values for the rows that already have a decision.
$ws.Cells.Item(3, 3).Value = "Đã chốt"
$ws.Cells.Item(4, 3).Value = "Đã chốt"
$ws.Cells.Item(6, 3).Value = "Đang trao đổi"

# Center the data cells of the new column (rows 2-16, under the header).
$range = $ws.Range("C2:C16")
$range.HorizontalAlignment = -4108  # xlCenter
$range.VerticalAlignment = -4108    # xlCenter

# Give the new column roughly the same width as column B.
$ws.Columns("C").ColumnWidth = 29.67

# Touch the column's default (below-the-data) formatting so the
# column picks up a center-aligned default style, then drop the
# scratch row again so the used range stays A1:G16.
$ws.Cells.Item(17, 3).HorizontalAlignment = -4108
$ws.Rows("17:17").Delete()

# Move the active selection, matching where the author left off editing.
$ws.Range("C7").Select()
